$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "69.362.41"
$ws.Range("E2").Value = "  +0.15%  "
$ws.Range("D3").Value = "3.419.60"
$ws.Range("E3").Value = "  +1.07%  "
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "582.09"
$ws.Range("E5").Value = "  -1.12%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "176.28"
$ws.Range("E6").Value = "  -2.40%  "
$ws.Range("E7").Value = "  +0.10%  "
$ws.Range("D8").Value = "3.411.42"
$ws.Range("E8").Value = "  +0.94%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.591"
$ws.Range("E9").Value = "  -0.58%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.198"
$ws.Range("E10").Value = "  +0.86%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.581"
$ws.Range("E11").Value = "  -1.47%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "48.55"
$ws.Range("E12").Value = "  -0.24%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000279"
$ws.Range("E13").Value = "  -2.63%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "690.95"
$ws.Range("E14").Value = "  +0.56%  "
$ws.Range("D15").Value = "3.971.90"
$ws.Range("E15").Value = "  +0.80%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "8.60"
$ws.Range("E16").Value = "  -0.34%  "
$ws.Range("D17").Value = "69.487.27"
$ws.Range("E17").Value = "  +0.14%  "
$ws.Range("D18").Value = "3.418.96"
$ws.Range("E18").Value = "  +2.27%  "
$ws.Range("E19").Value = "  +0.82%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "17.61"
$ws.Range("E20").Value = "  -0.82%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "11.36"
$ws.Range("E21").Value = "  -0.44%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.892"
$ws.Range("E22").Value = "  -1.07%  "
$ws.Range("E23").Value = "  +0.22%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "16.85"
$ws.Range("E24").Value = "  -1.44%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "100.60"
$ws.Range("E25").Value = "  -3.74%  "
$ws.Range("E26").Value = "  -1.43%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.66"
$ws.Range("E27").Value = "  -2.40%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.53"
$ws.Range("E28").Value = "  -1.00%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "33.40"
$ws.Range("E29").Value = "  -3.02%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "8.71"
$ws.Range("E30").Value = "  +0.10%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.00"
$ws.Range("E31").Value = "  +0.49%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "577.83"
$ws.Range("E32").Value = "  +3.64%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.67"
$ws.Range("E33").Value = "  +0.56%  "
$ws.Range("E34").Value = "  -1.96%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "58.37"
$ws.Range("E35").Value = "  +0.45%  "
$ws.Range("E36").Value = "  -3.43%  "
$ws.Range("D38").Value = "3.574.93"
$ws.Range("E38").Value = "  -3.50%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.138"
$ws.Range("E39").Value = "  -1.01%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "34.75"
$ws.Range("E40").Value = "  -0.57%  "
$ws.Range("E41").Value = "  +2.72%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.24"
$ws.Range("E42").Value = "  -0.36%  "
$ws.Range("E43").Value = "  -1.62%  "
$ws.Range("B44").Value = "ApeXProtocol"
$ws.Range("C44").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "3.36"
$ws.Range("E44").Value = "  +2.27%  "
$ws.Range("B45").Value = "TheGraph"
$ws.Range("C45").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.331"
$ws.Range("E45").Value = "  -2.57%  "
$ws.Range("B46").Value = "VeChain"
$ws.Range("C46").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0416"
$ws.Range("E46").Value = "  -0.77%  "
$ws.Range("B47").Value = "Mantle"
$ws.Range("C47").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.45"
$ws.Range("E47").Value = "  +4.11%  "
$ws.Range("B48").Value = "ThetaToken"
$ws.Range("C48").Value = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.64"
$ws.Range("E48").Value = "  -0.66%  "
$ws.Range("B49").Value = "Stellar"
$ws.Range("C49").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.128"
$ws.Range("E49").Value = "  -1.50%  "
$ws.Range("B50").Value = "FirstDigitalUSD"
$ws.Range("C50").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.00"
$ws.Range("E50").Value = "  -0.26%  "
$ws.Range("B51").Value = "Monero"
$ws.Range("C51").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "132.41"
$ws.Range("E51").Value = "  -0.48%  "
